$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Policy Number" column (A) is being dropped; "Source URL" (was B)
# moves into A and "Destination URL" (was C) moves into B. Shift the
# values left one column instead of doing a structural column delete so
# the existing column-width metadata for B/C carries over cleanly.
$lastRow = 11
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 3).Clear()
}

# Column A now holds URLs like column B used to, so give it the same width.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
